# "theorie kleine fouten verbeterd"
#
# The sentence "... de eerste is 1 waard, de 2e ..." gets the word
# "plek" inserted after "eerste" -> "... de eerste plek is 1 waard, de 2e ...".
# The (hidden) "_GoBack" last-edit bookmark -- which sat at the very end of
# the paragraph/document before the edit -- moves to sit right after the
# newly typed text ("de eerste plek "), reflecting where the author's
# cursor was left after making this edit.

$d = $word.ActiveDocument

# --- 1. Insert the missing word "plek" --------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute(
    "de eerste is 1 waard, de 2", $true, $false, $false, $false, $false,
    $true, 1, $false, "de eerste plek is 1 waard, de 2", 2)

if (-not $found) {
    throw "Could not find target sentence to fix."
}

# --- 2. Relocate the '_GoBack' bookmark --------------------------------
# Remove it from wherever it currently is (end of the last list item).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Re-create it right after "de eerste plek " (i.e. right before
# "is 1 waard, de 2"), matching where the author's last edit happened.
$markRng = $d.Content
$markRng.Find.ClearFormatting()
$markFound = $markRng.Find.Execute(
    "de eerste plek ", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)

if ($markFound) {
    $markRng.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $markRng)
}
